$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "normal attack" / "skill" label strings in column A (ID column) to
# their new upper-case forms, and set the "NextLevelID" column (C) so every
# row now points at NORMALATTACK2.
$ws.Range("A2").Value = "NORMALATTACK1"
$ws.Range("A3").Value = "NORMALATTACK2"
$ws.Range("A4").Value = "NORMALATTACK3"
$ws.Range("A5").Value = "NORMALTHUMP"
$ws.Range("A6").Value = "SKILL1"
$ws.Range("A7").Value = "SKILL2"
$ws.Range("A8").Value = "SKILL3"
$ws.Range("A9").Value = "SKILL4"

$ws.Range("C2").Value = "NORMALATTACK2"
$ws.Range("C3").Value = "NORMALATTACK2"
$ws.Range("C4").Value = "NORMALATTACK2"
$ws.Range("C5").Value = "NORMALATTACK2"
$ws.Range("C6").Value = "NORMALATTACK2"
$ws.Range("C7").Value = "NORMALATTACK2"
$ws.Range("C8").Value = "NORMALATTACK2"
$ws.Range("C9").Value = "NORMALATTACK2"

# Bump the AnimaState ids for the new "skill" rows.
$ws.Range("G6").Value = 101
$ws.Range("G7").Value = 102
$ws.Range("G8").Value = 103
$ws.Range("G9").Value = 104

# The NextLevelID column for rows 4-9 loses its shaded/bordered "new skill"
# look now that it is simply pointing at another normal attack, and is
# formatted as plain Text.
$ws.Range("C4:C9").ClearFormats()
$ws.Range("C4:C9").NumberFormat = "@"

# Mirror the author's final selection (two blocks of the NextLevelID column).
$ws.Range("C2:C5").Select() | Out-Null
$ws.Range("C6:C9").Select() | Out-Null
